# Commit: "Sun, May 03, 2020  1:05:49 AM"
#
# The canonical OOXML diff swaps the contents of ppt/theme/theme1.xml
# (the deck's "Integral" theme, used by the slide master) and
# ppt/theme/theme2.xml (the default "Office Theme", used by the notes
# master): the slide master ends up with the Office Theme's color
# palette and the notes master ends up with the old Integral palette.
#
# The font scheme and format scheme (fills/lines/effects) are already
# byte-for-byte identical between the two themes, so the only
# substantive, user-visible change is the 12-slot theme color palette.
# We rewrite that palette -- reachable through the presentation's
# theme-color-scheme API, which edits the deck's shared theme part --
# from the "Integral" colors to the standard "Office Theme" colors.

# Helper: build the VBA-style RGB() long (R + G*256 + B*65536) that the
# ThemeColorScheme RGBColor.RGB setter expects, from a "RRGGBB" hex string.
function HexToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette, in ThemeColorScheme index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeTheme = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeTheme.Count; $i++) {
    $tcs.Item($i + 1).RGB = HexToRgbLong $officeTheme[$i]
}
